# Atualizei dados add e bibi
# Atualiza a planilha de vendas atipicas: corrige estoque/desvio de algumas
# linhas existentes e adiciona as vendas do dia 2025-07-02.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrige valores recalculados em linhas ja existentes ---------------

# Linhas do produto "FONE BLUETOOTH BASIKE TWS FON6694" (id_produto 13079)
$ws.Range("G2").Value = -453
$ws.Range("I2").Value = 0.28

$ws.Range("G3").Value = -453
$ws.Range("I3").Value = 0.28

$ws.Range("G5").Value = -453
$ws.Range("I5").Value = 0.28

# Linha da "BALANCA DIGITAL 10KG" (id_produto 14186)
$ws.Range("G7").Value = -140

# Linha do "FONE SEM FIO BOX ... LETRON" em 2025-07-01 (id_venda 375082)
$ws.Range("G9").Value = -317
$ws.Range("I9").Value = 0.3

# --- Adiciona as novas vendas atipicas de 2025-07-02 ---------------------

# Linha 10: venda 375697 - ADAPTADOR TUDO EM UM UNIVERSAL INOVA PRIME TRA-30078
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-07-02"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "BEMOL S/A"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "375697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = 13018
$ws.Range("F10").Value = "ADAPTADOR TUDO EM UM UNIVERSAL INOVA PRIME TRA-30078"
$ws.Range("G10").Value = -6
$ws.Range("H10").Value = 1.08
$ws.Range("I10").Value = 0.29

# Linha 11: venda 375731 - FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025-07-02"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "BEMOL S/A"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "375731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = 13546
$ws.Range("F11").Value = "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Range("G11").Value = -317
$ws.Range("H11").Value = 1.09
$ws.Range("I11").Value = 0.3
